{"js": "const replacements = [\n  [\"11\u00d799=1089\", \"82\u00d728=2296\"],\n  [\"22\u00d711=242\", \"58\u00d711=638\"],\n  [\"30\u00d725=750\", \"93\u00d711=1023\"],\n  [\"48\u00d741=1968\", \"29\u00d751=1479\"],\n  [\"22\u00d733=726\", \"23\u00d764=1472\"],\n  [\"70\u00d759=4130\", \"90\u00d754=4860\"],\n  [\"25\u00d761=1525\", \"21\u00d727=567\"],\n  [\"79\u00d778=6162\", \"81\u00d788=7128\"],\n  [\"95\u00d799=9405\", \"47\u00d729=1363\"],\n  [\"45\u00d715=675\", \"22\u00d792=2024\"],\n  [\"67\u00d727=1809\", \"24\u00d746=1104\"],\n  [\"23\u00d743=989\", \"80\u00d726=2080\"],\n  [\"94\u00d798=9212\", \"64\u00d770=4480\"],\n  [\"34\u00d737=1258\", \"74\u00d718=1332\"],\n  [\"93\u00d781=7533\", \"62\u00d735=2170\"],\n  [\"78\u00d738=2964\", \"18\u00d737=666\"],\n  [\"89\u00d735=3115\", \"24\u00d735=840\"],\n  [\"25\u00d724=600\", \"89\u00d763=5607\"],\n  [\"76\u00d715=1140\", \"55\u00d756=3080\"],\n  [\"97\u00d755=5335\", \"41\u00d766=2706\"],\n  [\"43\u00d724=1032\", \"36\u00d740=1440\"],\n  [\"85\u00d791=7735\", \"61\u00d775=4575\"],\n  [\"16\u00d796=1536\", \"41\u00d771=2911\"],\n  [\"63\u00d746=2898\", \"22\u00d753=1166\"],\n  [\"43\u00d777=3311\", \"96\u00d799=9504\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @{\n    \"11\u00d799=1089\" = \"82\u00d728=2296\"\n    \"22\u00d711=242\"  = \"58\u00d711=638\"\n    \"30\u00d725=750\"  = \"93\u00d711=1023\"\n    \"48\u00d741=1968\" = \"29\u00d751=1479\"\n    \"22\u00d733=726\"  = \"23\u00d764=1472\"\n    \"70\u00d759=4130\" = \"90\u00d754=4860\"\n    \"25\u00d761=1525\" = \"21\u00d727=567\"\n    \"79\u00d778=6162\" = \"81\u00d788=7128\"\n    \"95\u00d799=9405\" = \"47\u00d729=1363\"\n    \"45\u00d715=675\"  = \"22\u00d792=2024\"\n    \"67\u00d727=1809\" = \"24\u00d746=1104\"\n    \"23\u00d743=989\"  = \"80\u00d726=2080\"\n    \"94\u00d798=9212\" = \"64\u00d770=4480\"\n    \"34\u00d737=1258\" = \"74\u00d718=1332\"\n    \"93\u00d781=7533\" = \"62\u00d735=2170\"\n    \"78\u00d738=2964\" = \"18\u00d737=666\"\n    \"89\u00d735=3115\" = \"24\u00d735=840\"\n    \"25\u00d724=600\"  = \"89\u00d763=5607\"\n    \"76\u00d715=1140\" = \"55\u00d756=3080\"\n    \"97\u00d755=5335\" = \"41\u00d766=2706\"\n    \"43\u00d724=1032\" = \"36\u00d740=1440\"\n    \"85\u00d791=7735\" = \"61\u00d775=4575\"\n    \"16\u00d796=1536\" = \"41\u00d771=2911\"\n    \"63\u00d746=2898\" = \"22\u00d753=1166\"\n    \"43\u00d777=3311\" = \"96\u00d799=9504\"\n}\n\nforeach ($old in $replacements.Keys) {\n    $new = $replacements[$old]\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
